$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REGION1-TSData")

# Insert a new row above the current row 2 (COM_CSTNET row), shifting it down to row 3
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with the COM_BNDNET data
$ws.Range("B2").Value = "COM_BNDNET"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "CO2EQS"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "ANNUAL"
$ws.Range("H2").Value = "UP"
$ws.Range("I2").Value = 0

$nullCols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP")
foreach ($col in $nullCols) {
    $ws.Range($col + "2").Value = "null"
}

$ws.Range("AB2").Value = 10000

$wb.Save()
